$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "64.857.09"
$ws.Range("E2").Value = "  -0.56%  "

$ws.Range("D3").Value = "3.451.53"
$ws.Range("E3").Value = "  -0.92%  "

$ws.Range("E4").Value = "  -0.04%  "

Set-TextValue $ws.Range("D5") "573.24"
$ws.Range("E5").Value = "  -1.04%  "

Set-TextValue $ws.Range("D6") "159.11"
$ws.Range("E6").Value = "  -1.85%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "3.446.58"
$ws.Range("E8").Value = "  -1.05%  "

Set-TextValue $ws.Range("D9") "0.575"
$ws.Range("E9").Value = "  -6.51%  "

Set-TextValue $ws.Range("D10") "7.19"
$ws.Range("E10").Value = "  -1.01%  "

$ws.Range("E11").Value = "  -3.53%  "

Set-TextValue $ws.Range("D12") "0.441"
$ws.Range("E12").Value = "  -1.21%  "

$ws.Range("D13").Value = "4.038.52"
$ws.Range("E13").Value = "  -1.22%  "

$ws.Range("E14").Value = "  -0.42%  "

Set-TextValue $ws.Range("D15") "27.68"
$ws.Range("E15").Value = "  -3.61%  "

$ws.Range("E16").Value = "  -9.80%  "

$ws.Range("D17").Value = "64.893.35"
$ws.Range("E17").Value = "  -0.50%  "

$ws.Range("D18").Value = "3.446.76"
$ws.Range("E18").Value = "  -0.95%  "

Set-TextValue $ws.Range("D19") "6.24"
$ws.Range("E19").Value = "  -3.26%  "

Set-TextValue $ws.Range("D20") "13.73"
$ws.Range("E20").Value = "  -4.43%  "

Set-TextValue $ws.Range("D21") "378.10"
$ws.Range("E21").Value = "  -1.26%  "

Set-TextValue $ws.Range("D22") "7.98"
$ws.Range("E22").Value = "  -3.00%  "

$ws.Range("E23").Value = "  +0.27%  "

Set-TextValue $ws.Range("D24") "0.542"
$ws.Range("E24").Value = "  -1.73%  "

Set-TextValue $ws.Range("D25") "72.05"
$ws.Range("E25").Value = "  -1.02%  "

$ws.Range("E26").Value = "  -0.28%  "

Set-TextValue $ws.Range("D27") "9.85"
$ws.Range("E27").Value = "  -2.62%  "

Set-TextValue $ws.Range("D28") "0.179"
$ws.Range("E28").Value = "  -0.31%  "

Set-TextValue $ws.Range("D29") "1.00"
$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("E30").Value = "  -4.67%  "

Set-TextValue $ws.Range("D31") "6.07"
$ws.Range("E31").Value = "  -2.75%  "

$ws.Range("E32").Value = "  -2.54%  "

Set-TextValue $ws.Range("D33") "23.17"
$ws.Range("E33").Value = "  -2.07%  "

Set-TextValue $ws.Range("D34") "6.97"
$ws.Range("E34").Value = "  -3.40%  "

Set-TextValue $ws.Range("D35") "1.56"
$ws.Range("E35").Value = "  -2.19%  "

Set-TextValue $ws.Range("D36") "160.95"
$ws.Range("E36").Value = "  -0.56%  "

Set-TextValue $ws.Range("D37") "1.87"
$ws.Range("E37").Value = "  -2.92%  "

$ws.Range("D38").Value = "2.914.97"
$ws.Range("E38").Value = "  -4.55%  "

Set-TextValue $ws.Range("D39") "0.0747"
$ws.Range("E39").Value = "  -3.55%  "

Set-TextValue $ws.Range("D40") "6.69"
$ws.Range("E40").Value = "  +1.79%  "

Set-TextValue $ws.Range("D41") "26.18"
$ws.Range("E41").Value = "  -3.25%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D42") "4.54"
$ws.Range("E42").Value = "  -0.88%  "

$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D43") "42.98"
$ws.Range("E43").Value = "  +0.18%  "

Set-TextValue $ws.Range("D44") "0.783"
$ws.Range("E44").Value = "  +0.52%  "

Set-TextValue $ws.Range("D45") "26.02"
$ws.Range("E45").Value = "  +0.15%  "

$ws.Range("E46").Value = "  -3.32%  "

Set-TextValue $ws.Range("D47") "2.32"
$ws.Range("E47").Value = "  +5.06%  "

Set-TextValue $ws.Range("D48") "319.17"
$ws.Range("E48").Value = "  +0.36%  "

$ws.Range("E49").Value = "  -2.63%  "

Set-TextValue $ws.Range("D50") "6.47"
$ws.Range("E50").Value = "  -4.28%  "

Set-TextValue $ws.Range("D51") "0.848"
$ws.Range("E51").Value = "  -3.36%  "
